# Applies the "update performer after changes from uat" edit to Config.xlsx
# Target sheets: Settings (sheet1), Assets (sheet3). Constants (sheet2) is untouched.

$wb = $excel.ActiveWorkbook
$wsSettings = $wb.Worksheets.Item("Settings")
$wsAssets   = $wb.Worksheets.Item("Assets")

# ---------------------------------------------------------------------------
# 1. Settings sheet: simple in-place text updates (rows 1-33 keep their shape)
# ---------------------------------------------------------------------------
$wsSettings.Range("B5").Value  = "Recharges Vodafone Usage Performer"
$wsSettings.Range("A18").Value = "HeaderSheetName"

# ---------------------------------------------------------------------------
# 2. Settings sheet: insert 9 blank rows before the old "Active List" header
#    (old row 35) so that everything from there on shifts down by 9 rows,
#    and fill rows 34-43 with the new "mapping / header" block of settings.
# ---------------------------------------------------------------------------
$wsSettings.Rows("35:43").Insert()

$newBlock = @(
    @(34, "VodafoneMappingSheetName", "Active Names"),
    @(35, "PreparedByRange",          "F15"),
    @(36, "DatePreparedRange",        "H15"),
    @(37, "DateAuthorizedRange",      "H16"),
    @(38, "DatePostedRange",          "H17"),
    @(39, "FinancialYearRange",       "B9"),
    @(40, "DocumentDateRange",        "B13"),
    @(41, "PostingDateRange",         "B14"),
    @(42, "PostingPeriodRange",       "B16"),
    @(43, "PreparedByName",           "Robot")
)
foreach ($item in $newBlock) {
    $r = $item[0]
    $wsSettings.Range("A$r").Value = $item[1]
    $wsSettings.Range("B$r").Value = $item[2]
}

# ---------------------------------------------------------------------------
# 3. Settings sheet: update the (now shifted) Mail section entry and append
#    a new row to the VBA functions section.
# ---------------------------------------------------------------------------
$wsSettings.Range("A52").Value = "MailSubject_Journal"
$wsSettings.Range("B52").Value = "Recharges Automation - Vodafone Usage - Journal {0}"

$wsSettings.Range("A58").Value = "vbaRemoveDuplicatesInColumn"
$wsSettings.Range("B58").Value = "RemoveDuplicatesInColumn"

# ---------------------------------------------------------------------------
# 4. Settings sheet: extend the trailing empty rows down to row 1025
#    (mirrors the extra blank formatted rows added at the bottom of sheet).
# ---------------------------------------------------------------------------
$wsSettings.Rows("1017:1025").RowHeight = 14.25

# ---------------------------------------------------------------------------
# 5. Assets sheet updates
# ---------------------------------------------------------------------------
$wsAssets.Range("A3").Value = "Recharges_Mail_Account"

$wsAssets.Range("A4").Value = "Asset_MAIL_TO"
$wsAssets.Range("B4").Value = "Recharges_MAIL_TO"

$wsAssets.Range("A5").Value = "Recharges_VodafoneUsageMappingFile"
$wsAssets.Range("B5").Value = "Recharges_VodafoneUsageMappingFile"

# ---------------------------------------------------------------------------
# 6. View state: Assets selection moves to A4, Settings selection moves to
#    A43 (scrolled further down); Settings stays the active/selected sheet.
# ---------------------------------------------------------------------------
$wsAssets.Activate()
$wsAssets.Range("A4").Select()

$wsSettings.Activate()
$wsSettings.Range("A43").Select()
